$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append three new data rows (8, 9, 10) below the existing table (A1:E7),
# mirroring the pattern of the existing rows. The ID/Projeto/Data/Local
# columns hold numeric- and date-looking text that must stay as literal
# text (not be coerced into a number/date), so each value is entered with
# a leading apostrophe (quote-prefix) the way a user typing into the grid
# would force text entry; the leftover quote-prefix formatting is then
# cleared so the cells end up with plain default formatting.

$rows = @(
    @{ Row = 8;  A = "-421549506"; B = "152262892521"; C = "20/12/2024"; D = "kjhhhuhhiu"; E = "wusgdhr" },
    @{ Row = 9;  A = "-421549506"; B = "152262892521"; C = "20/12/2024"; D = "kjhhhuhhiu"; E = "N/A" },
    @{ Row = 10; A = "-421549506"; B = "152262892521"; C = "20/12/2024"; D = "kjhhhuhhiu"; E = "wusgdhri" }
)

foreach ($r in $rows) {
    $n = $r.Row
    $ws.Range("A$n").Value = "'" + $r.A
    $ws.Range("B$n").Value = "'" + $r.B
    $ws.Range("C$n").Value = "'" + $r.C
    $ws.Range("D$n").Value = $r.D
    $ws.Range("E$n").Value = $r.E
}

# Drop the quote-prefix formatting picked up above so the new cells match
# the plain/default styling of the rest of the sheet.
$ws.Range("A8:E10").ClearFormats()
